$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where week 6 (column H) attendance should be marked as present (1)
$rows = @(6, 9, 12, 13, 14, 18, 19, 21)

foreach ($r in $rows) {
    $srcCell = $ws.Cells.Item($r, 7)   # column G, same row, used as style template
    $dstCell = $ws.Cells.Item($r, 8)   # column H
    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $dstCell.Value = 1
}

$excel.CutCopyMode = 0

# Update the active selection on the sheet to match the saved view state
$ws.Range("H13").Select() | Out-Null
